# Merge 'TSM Cord of the Beresaad'
# Adds a new tlk string entry (id 6610058 / 6610059) for the amulet
# "Cord of the Beresaad" plus its description, and annotates the new
# entry with the same "TSM <name>" authoring comment used by the other
# entries in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow1 = 60
$newRow2 = 61

# --- new tlk rows -----------------------------------------------------
$ws.Cells.Item($newRow1, 1).Value = 6610058
$ws.Cells.Item($newRow1, 2).Value = "Cord of the Beresaad"

$ws.Cells.Item($newRow2, 1).Value = 6610059
$ws.Cells.Item($newRow2, 2).Value = "Intricately braided of many thin leather strands, this sturdy torc is otherwise unornamented."

# --- match formatting of the surrounding rows --------------------------
# Column A / B on the existing rows carry a custom (unnamed) cell style;
# copy it across instead of recreating the xf by hand.
$ws.Range("A2").Copy()
$ws.Range("A60:A61").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B60:B61").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- annotate the new entry, like the other "TSM ..." rows -------------
$ws.Range("A60").AddComment("TSM Cord of the Beresaad")

# --- mirror the selection state recorded in the saved workbook ---------
[void]$ws.Range("B74").Select()
